# =====================================================================
# Commit: "feat: add 2022-Q3 data"
#
# 1) Insert a new worksheet "2022-Q3" right after "总计" (shifting the
#    older quarterly sheets down by one tab position; their own data is
#    untouched).
# 2) On "总计", add a new summary row for 2022-Q3 (29 holdings, 28.44
#    billion yuan), pushing the existing history rows down by one.
# 3) Fill the new "2022-Q3" sheet with the fund holdings table.
#
# NOTE: reading `Range.Value` back out as an rvalue is unreliable in
# this host (it can yield the property descriptor instead of the cell
# contents), so every value below is written as a literal instead of
# being copied from another cell at runtime.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("总计")

# ---- Step 1: insert new worksheet right after "总计" ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2022-Q3"

# ---- Step 2: rewrite "总计" rows 2-9 (row 1 header is unchanged) ----
$styleA = $ws1.Range("A2")  # style index 2 (bold/border) already applied here

# row 2: 2022-Q3
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 29
$ws1.Range("D2").Value = 28.44

# row 3: 2022-Q2
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 35
$ws1.Range("D3").Value = 48.82

# row 4: 2022-Q1
$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 37
$ws1.Range("D4").Value = 51.24

# row 5: 2021-Q4
$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 41
$ws1.Range("D5").Value = 54.05

# row 6: 2021-Q3
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 33
$ws1.Range("D6").Value = 52.24

# row 7: 2021-Q2
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "2021-Q2"
$ws1.Range("C7").Value = 83
$ws1.Range("D7").Value = 77.55

# row 8: 2021-Q1
$ws1.Range("A8").Value = 6
$ws1.Range("B8").Value = "2021-Q1"
$ws1.Range("C8").Value = 63
$ws1.Range("D8").Value = 64.18

# row 9: 2020-Q4
$cell = $ws1.Range("A9"); $styleA.Copy($cell); $cell.Value = 7
$ws1.Range("B9").Value = "2020-Q4"
$ws1.Range("C9").Value = 80
$ws1.Range("D9").Value = 71.84

# ---- Step 3: populate "2022-Q3" sheet ----
$styleSrc = $ws1.Range("A2")  # reuse the bold/bordered "total" style (style index 2)

# Header row
$cell = $ws2.Range("B1"); $styleSrc.Copy($cell); $cell.Value = "基金代码"
$cell = $ws2.Range("C1"); $styleSrc.Copy($cell); $cell.Value = "基金名称"
$cell = $ws2.Range("D1"); $styleSrc.Copy($cell); $cell.Value = "基金规模"
$cell = $ws2.Range("E1"); $styleSrc.Copy($cell); $cell.Value = "股票总仓位"
$cell = $ws2.Range("F1"); $styleSrc.Copy($cell); $cell.Value = "仓位占比"
$cell = $ws2.Range("G1"); $styleSrc.Copy($cell); $cell.Value = "持有市值(亿元)"
$cell = $ws2.Range("H1"); $styleSrc.Copy($cell); $cell.Value = "仓位排名"

# row 2: index 0 -> 000751 嘉实新兴产业股票
$cell = $ws2.Range("A2"); $styleSrc.Copy($cell); $cell.Value = 0
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "000751"
$ws2.Range("C2").Value = "嘉实新兴产业股票"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "76.18"
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "92.94"
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "7.01"
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "5.3402"
$ws2.Range("H2").Value = 4

# row 3: index 1 -> 010186 嘉实核心成长混合A
$cell = $ws2.Range("A3"); $styleSrc.Copy($cell); $cell.Value = 1
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "010186"
$ws2.Range("C3").Value = "嘉实核心成长混合A"
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "60.62"
$ws2.Range("E3").NumberFormat = "@"
$ws2.Range("E3").Value = "91.81"
$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Value = "7.42"
$ws2.Range("G3").NumberFormat = "@"
$ws2.Range("G3").Value = "4.4980"
$ws2.Range("H3").Value = 3

# row 4: index 2 -> 009795 嘉实远见精选两年持有期混合
$cell = $ws2.Range("A4"); $styleSrc.Copy($cell); $cell.Value = 2
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "009795"
$ws2.Range("C4").Value = "嘉实远见精选两年持有期混合"
$ws2.Range("D4").NumberFormat = "@"
$ws2.Range("D4").Value = "55.77"
$ws2.Range("E4").NumberFormat = "@"
$ws2.Range("E4").Value = "93.56"
$ws2.Range("F4").NumberFormat = "@"
$ws2.Range("F4").Value = "7.29"
$ws2.Range("G4").NumberFormat = "@"
$ws2.Range("G4").Value = "4.0656"
$ws2.Range("H4").Value = 4

# row 5: index 3 -> 000595 嘉实泰和混合
$cell = $ws2.Range("A5"); $styleSrc.Copy($cell); $cell.Value = 3
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "000595"
$ws2.Range("C5").Value = "嘉实泰和混合"
$ws2.Range("D5").NumberFormat = "@"
$ws2.Range("D5").Value = "40.15"
$ws2.Range("E5").NumberFormat = "@"
$ws2.Range("E5").Value = "92.39"
$ws2.Range("F5").NumberFormat = "@"
$ws2.Range("F5").Value = "6.86"
$ws2.Range("G5").NumberFormat = "@"
$ws2.Range("G5").Value = "2.7543"
$ws2.Range("H5").Value = 4

# row 6: index 4 -> 000011 华夏大盘精选混合A
$cell = $ws2.Range("A6"); $styleSrc.Copy($cell); $cell.Value = 4
$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = "000011"
$ws2.Range("C6").Value = "华夏大盘精选混合A"
$ws2.Range("D6").NumberFormat = "@"
$ws2.Range("D6").Value = "40.09"
$ws2.Range("E6").NumberFormat = "@"
$ws2.Range("E6").Value = "89.60"
$ws2.Range("F6").NumberFormat = "@"
$ws2.Range("F6").Value = "4.46"
$ws2.Range("G6").NumberFormat = "@"
$ws2.Range("G6").Value = "1.7880"
$ws2.Range("H6").Value = 9

# row 7: index 5 -> 070002 嘉实增长混合
$cell = $ws2.Range("A7"); $styleSrc.Copy($cell); $cell.Value = 5
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = "070002"
$ws2.Range("C7").Value = "嘉实增长混合"
$ws2.Range("D7").NumberFormat = "@"
$ws2.Range("D7").Value = "27.24"
$ws2.Range("E7").NumberFormat = "@"
$ws2.Range("E7").Value = "74.18"
$ws2.Range("F7").NumberFormat = "@"
$ws2.Range("F7").Value = "5.91"
$ws2.Range("G7").NumberFormat = "@"
$ws2.Range("G7").Value = "1.6099"
$ws2.Range("H7").Value = 5

# row 8: index 6 -> 100026 富国天合稳健混合
$cell = $ws2.Range("A8"); $styleSrc.Copy($cell); $cell.Value = 6
$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "100026"
$ws2.Range("C8").Value = "富国天合稳健混合"
$ws2.Range("D8").NumberFormat = "@"
$ws2.Range("D8").Value = "41.01"
$ws2.Range("E8").NumberFormat = "@"
$ws2.Range("E8").Value = "78.83"
$ws2.Range("F8").NumberFormat = "@"
$ws2.Range("F8").Value = "3.44"
$ws2.Range("G8").NumberFormat = "@"
$ws2.Range("G8").Value = "1.4107"
$ws2.Range("H8").Value = 8

# row 9: index 7 -> 166025 中欧远见两年定期开放混合A
$cell = $ws2.Range("A9"); $styleSrc.Copy($cell); $cell.Value = 7
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "166025"
$ws2.Range("C9").Value = "中欧远见两年定期开放混合A"
$ws2.Range("D9").NumberFormat = "@"
$ws2.Range("D9").Value = "44.75"
$ws2.Range("E9").NumberFormat = "@"
$ws2.Range("E9").Value = "59.87"
$ws2.Range("F9").NumberFormat = "@"
$ws2.Range("F9").Value = "2.93"
$ws2.Range("G9").NumberFormat = "@"
$ws2.Range("G9").Value = "1.3112"
$ws2.Range("H9").Value = 8

# row 10: index 8 -> 009137 嘉实瑞和两年持有期混合
$cell = $ws2.Range("A10"); $styleSrc.Copy($cell); $cell.Value = 8
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "009137"
$ws2.Range("C10").Value = "嘉实瑞和两年持有期混合"
$ws2.Range("D10").NumberFormat = "@"
$ws2.Range("D10").Value = "17.33"
$ws2.Range("E10").NumberFormat = "@"
$ws2.Range("E10").Value = "93.36"
$ws2.Range("F10").NumberFormat = "@"
$ws2.Range("F10").Value = "7.36"
$ws2.Range("G10").NumberFormat = "@"
$ws2.Range("G10").Value = "1.2755"
$ws2.Range("H10").Value = 5

# row 11: index 9 -> 002593 富国美丽中国混合A
$cell = $ws2.Range("A11"); $styleSrc.Copy($cell); $cell.Value = 9
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "002593"
$ws2.Range("C11").Value = "富国美丽中国混合A"
$ws2.Range("D11").NumberFormat = "@"
$ws2.Range("D11").Value = "28.84"
$ws2.Range("E11").NumberFormat = "@"
$ws2.Range("E11").Value = "82.00"
$ws2.Range("F11").NumberFormat = "@"
$ws2.Range("F11").Value = "3.21"
$ws2.Range("G11").NumberFormat = "@"
$ws2.Range("G11").Value = "0.9258"
$ws2.Range("H11").Value = 9

# row 12: index 10 -> 360006 光大保德信新增长混合
$cell = $ws2.Range("A12"); $styleSrc.Copy($cell); $cell.Value = 10
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "360006"
$ws2.Range("C12").Value = "光大保德信新增长混合"
$ws2.Range("D12").NumberFormat = "@"
$ws2.Range("D12").Value = "17.68"
$ws2.Range("E12").NumberFormat = "@"
$ws2.Range("E12").Value = "83.30"
$ws2.Range("F12").NumberFormat = "@"
$ws2.Range("F12").Value = "4.26"
$ws2.Range("G12").NumberFormat = "@"
$ws2.Range("G12").Value = "0.7532"
$ws2.Range("H12").Value = 5

# row 13: index 11 -> 519035 富国天博创新混合
$cell = $ws2.Range("A13"); $styleSrc.Copy($cell); $cell.Value = 11
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "519035"
$ws2.Range("C13").Value = "富国天博创新混合"
$ws2.Range("D13").NumberFormat = "@"
$ws2.Range("D13").Value = "18.16"
$ws2.Range("E13").NumberFormat = "@"
$ws2.Range("E13").Value = "91.30"
$ws2.Range("F13").NumberFormat = "@"
$ws2.Range("F13").Value = "3.07"
$ws2.Range("G13").NumberFormat = "@"
$ws2.Range("G13").Value = "0.5575"
$ws2.Range("H13").Value = 7

# row 14: index 12 -> 010187 嘉实核心成长混合C
$cell = $ws2.Range("A14"); $styleSrc.Copy($cell); $cell.Value = 12
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "010187"
$ws2.Range("C14").Value = "嘉实核心成长混合C"
$ws2.Range("D14").NumberFormat = "@"
$ws2.Range("D14").Value = "4.13"
$ws2.Range("E14").NumberFormat = "@"
$ws2.Range("E14").Value = "91.81"
$ws2.Range("F14").NumberFormat = "@"
$ws2.Range("F14").Value = "7.42"
$ws2.Range("G14").NumberFormat = "@"
$ws2.Range("G14").Value = "0.3064"
$ws2.Range("H14").Value = 3

# row 15: index 13 -> 070022 嘉实领先成长混合
$cell = $ws2.Range("A15"); $styleSrc.Copy($cell); $cell.Value = 13
$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "070022"
$ws2.Range("C15").Value = "嘉实领先成长混合"
$ws2.Range("D15").NumberFormat = "@"
$ws2.Range("D15").Value = "4.99"
$ws2.Range("E15").NumberFormat = "@"
$ws2.Range("E15").Value = "86.27"
$ws2.Range("F15").NumberFormat = "@"
$ws2.Range("F15").Value = "5.65"
$ws2.Range("G15").NumberFormat = "@"
$ws2.Range("G15").Value = "0.2819"
$ws2.Range("H15").Value = 5

# row 16: index 14 -> 000513 富国高端制造行业股票A
$cell = $ws2.Range("A16"); $styleSrc.Copy($cell); $cell.Value = 14
$ws2.Range("B16").NumberFormat = "@"
$ws2.Range("B16").Value = "000513"
$ws2.Range("C16").Value = "富国高端制造行业股票A"
$ws2.Range("D16").NumberFormat = "@"
$ws2.Range("D16").Value = "7.01"
$ws2.Range("E16").NumberFormat = "@"
$ws2.Range("E16").Value = "91.41"
$ws2.Range("F16").NumberFormat = "@"
$ws2.Range("F16").Value = "3.12"
$ws2.Range("G16").NumberFormat = "@"
$ws2.Range("G16").Value = "0.2187"
$ws2.Range("H16").Value = 9

# row 17: index 15 -> 001759 嘉实成长增强灵活配置混合
$cell = $ws2.Range("A17"); $styleSrc.Copy($cell); $cell.Value = 15
$ws2.Range("B17").NumberFormat = "@"
$ws2.Range("B17").Value = "001759"
$ws2.Range("C17").Value = "嘉实成长增强灵活配置混合"
$ws2.Range("D17").NumberFormat = "@"
$ws2.Range("D17").Value = "4.25"
$ws2.Range("E17").NumberFormat = "@"
$ws2.Range("E17").Value = "90.77"
$ws2.Range("F17").NumberFormat = "@"
$ws2.Range("F17").Value = "4.81"
$ws2.Range("G17").NumberFormat = "@"
$ws2.Range("G17").Value = "0.2044"
$ws2.Range("H17").Value = 7

# row 18: index 16 -> 005241 中欧时代智慧混合A
$cell = $ws2.Range("A18"); $styleSrc.Copy($cell); $cell.Value = 16
$ws2.Range("B18").NumberFormat = "@"
$ws2.Range("B18").Value = "005241"
$ws2.Range("C18").Value = "中欧时代智慧混合A"
$ws2.Range("D18").NumberFormat = "@"
$ws2.Range("D18").Value = "6.95"
$ws2.Range("E18").NumberFormat = "@"
$ws2.Range("E18").Value = "77.65"
$ws2.Range("F18").NumberFormat = "@"
$ws2.Range("F18").Value = "2.89"
$ws2.Range("G18").NumberFormat = "@"
$ws2.Range("G18").Value = "0.2009"
$ws2.Range("H18").Value = 10

# row 19: index 17 -> 011921 富国均衡成长三年持有期混合A
$cell = $ws2.Range("A19"); $styleSrc.Copy($cell); $cell.Value = 17
$ws2.Range("B19").NumberFormat = "@"
$ws2.Range("B19").Value = "011921"
$ws2.Range("C19").Value = "富国均衡成长三年持有期混合A"
$ws2.Range("D19").NumberFormat = "@"
$ws2.Range("D19").Value = "6.14"
$ws2.Range("E19").NumberFormat = "@"
$ws2.Range("E19").Value = "90.69"
$ws2.Range("F19").NumberFormat = "@"
$ws2.Range("F19").Value = "3.10"
$ws2.Range("G19").NumberFormat = "@"
$ws2.Range("G19").Value = "0.1903"
$ws2.Range("H19").Value = 8

# row 20: index 18 -> 001036 嘉实企业变革股票
$cell = $ws2.Range("A20"); $styleSrc.Copy($cell); $cell.Value = 18
$ws2.Range("B20").NumberFormat = "@"
$ws2.Range("B20").Value = "001036"
$ws2.Range("C20").Value = "嘉实企业变革股票"
$ws2.Range("D20").NumberFormat = "@"
$ws2.Range("D20").Value = "3.89"
$ws2.Range("E20").NumberFormat = "@"
$ws2.Range("E20").Value = "91.40"
$ws2.Range("F20").NumberFormat = "@"
$ws2.Range("F20").Value = "4.88"
$ws2.Range("G20").NumberFormat = "@"
$ws2.Range("G20").Value = "0.1898"
$ws2.Range("H20").Value = 7

# row 21: index 19 -> 512330 南方中证500信息技术指数ETF
$cell = $ws2.Range("A21"); $styleSrc.Copy($cell); $cell.Value = 19
$ws2.Range("B21").NumberFormat = "@"
$ws2.Range("B21").Value = "512330"
$ws2.Range("C21").Value = "南方中证500信息技术指数ETF"
$ws2.Range("D21").NumberFormat = "@"
$ws2.Range("D21").Value = "4.19"
$ws2.Range("E21").NumberFormat = "@"
$ws2.Range("E21").Value = "99.61"
$ws2.Range("F21").NumberFormat = "@"
$ws2.Range("F21").Value = "3.67"
$ws2.Range("G21").NumberFormat = "@"
$ws2.Range("G21").Value = "0.1538"
$ws2.Range("H21").Value = 4

# row 22: index 20 -> 008138 富国龙头优势混合
$cell = $ws2.Range("A22"); $styleSrc.Copy($cell); $cell.Value = 20
$ws2.Range("B22").NumberFormat = "@"
$ws2.Range("B22").Value = "008138"
$ws2.Range("C22").Value = "富国龙头优势混合"
$ws2.Range("D22").NumberFormat = "@"
$ws2.Range("D22").Value = "4.41"
$ws2.Range("E22").NumberFormat = "@"
$ws2.Range("E22").Value = "92.58"
$ws2.Range("F22").NumberFormat = "@"
$ws2.Range("F22").Value = "3.15"
$ws2.Range("G22").NumberFormat = "@"
$ws2.Range("G22").Value = "0.1389"
$ws2.Range("H22").Value = 9

# row 23: index 21 -> 007101 中欧远见两年定期开放混合C
$cell = $ws2.Range("A23"); $styleSrc.Copy($cell); $cell.Value = 21
$ws2.Range("B23").NumberFormat = "@"
$ws2.Range("B23").Value = "007101"
$ws2.Range("C23").Value = "中欧远见两年定期开放混合C"
$ws2.Range("D23").NumberFormat = "@"
$ws2.Range("D23").Value = "2.79"
$ws2.Range("E23").NumberFormat = "@"
$ws2.Range("E23").Value = "59.87"
$ws2.Range("F23").NumberFormat = "@"
$ws2.Range("F23").Value = "2.93"
$ws2.Range("G23").NumberFormat = "@"
$ws2.Range("G23").Value = "0.0817"
$ws2.Range("H23").Value = 8

# row 24: index 22 -> 003292 嘉实优势成长灵活配置混合
$cell = $ws2.Range("A24"); $styleSrc.Copy($cell); $cell.Value = 22
$ws2.Range("B24").NumberFormat = "@"
$ws2.Range("B24").Value = "003292"
$ws2.Range("C24").Value = "嘉实优势成长灵活配置混合"
$ws2.Range("D24").NumberFormat = "@"
$ws2.Range("D24").Value = "2.15"
$ws2.Range("E24").NumberFormat = "@"
$ws2.Range("E24").Value = "84.39"
$ws2.Range("F24").NumberFormat = "@"
$ws2.Range("F24").Value = "3.38"
$ws2.Range("G24").NumberFormat = "@"
$ws2.Range("G24").Value = "0.0727"
$ws2.Range("H24").Value = 7

# row 25: index 23 -> 005242 中欧时代智慧混合C
$cell = $ws2.Range("A25"); $styleSrc.Copy($cell); $cell.Value = 23
$ws2.Range("B25").NumberFormat = "@"
$ws2.Range("B25").Value = "005242"
$ws2.Range("C25").Value = "中欧时代智慧混合C"
$ws2.Range("D25").NumberFormat = "@"
$ws2.Range("D25").Value = "2.29"
$ws2.Range("E25").NumberFormat = "@"
$ws2.Range("E25").Value = "77.65"
$ws2.Range("F25").NumberFormat = "@"
$ws2.Range("F25").Value = "2.89"
$ws2.Range("G25").NumberFormat = "@"
$ws2.Range("G25").Value = "0.0662"
$ws2.Range("H25").Value = 10

# row 26: index 24 -> 011566 富国美丽中国混合C
$cell = $ws2.Range("A26"); $styleSrc.Copy($cell); $cell.Value = 24
$ws2.Range("B26").NumberFormat = "@"
$ws2.Range("B26").Value = "011566"
$ws2.Range("C26").Value = "富国美丽中国混合C"
$ws2.Range("D26").NumberFormat = "@"
$ws2.Range("D26").Value = "0.71"
$ws2.Range("E26").NumberFormat = "@"
$ws2.Range("E26").Value = "82.00"
$ws2.Range("F26").NumberFormat = "@"
$ws2.Range("F26").Value = "3.21"
$ws2.Range("G26").NumberFormat = "@"
$ws2.Range("G26").Value = "0.0228"
$ws2.Range("H26").Value = 9

# row 27: index 25 -> 011922 富国均衡成长三年持有期混合C
$cell = $ws2.Range("A27"); $styleSrc.Copy($cell); $cell.Value = 25
$ws2.Range("B27").NumberFormat = "@"
$ws2.Range("B27").Value = "011922"
$ws2.Range("C27").Value = "富国均衡成长三年持有期混合C"
$ws2.Range("D27").NumberFormat = "@"
$ws2.Range("D27").Value = "0.44"
$ws2.Range("E27").NumberFormat = "@"
$ws2.Range("E27").Value = "90.69"
$ws2.Range("F27").NumberFormat = "@"
$ws2.Range("F27").Value = "3.10"
$ws2.Range("G27").NumberFormat = "@"
$ws2.Range("G27").Value = "0.0136"
$ws2.Range("H27").Value = 8

# row 28: index 26 -> 012628 华夏大盘精选混合C
$cell = $ws2.Range("A28"); $styleSrc.Copy($cell); $cell.Value = 26
$ws2.Range("B28").NumberFormat = "@"
$ws2.Range("B28").Value = "012628"
$ws2.Range("C28").Value = "华夏大盘精选混合C"
$ws2.Range("D28").NumberFormat = "@"
$ws2.Range("D28").Value = "0.17"
$ws2.Range("E28").NumberFormat = "@"
$ws2.Range("E28").Value = "89.60"
$ws2.Range("F28").NumberFormat = "@"
$ws2.Range("F28").Value = "4.46"
$ws2.Range("G28").NumberFormat = "@"
$ws2.Range("G28").Value = "0.0076"
$ws2.Range("H28").Value = 9

# row 29: index 27 -> 561150 富国中证500ESG基准ETF
$cell = $ws2.Range("A29"); $styleSrc.Copy($cell); $cell.Value = 27
$ws2.Range("B29").NumberFormat = "@"
$ws2.Range("B29").Value = "561150"
$ws2.Range("C29").Value = "富国中证500ESG基准ETF"
$ws2.Range("D29").NumberFormat = "@"
$ws2.Range("D29").Value = "0.35"
$ws2.Range("E29").NumberFormat = "@"
$ws2.Range("E29").Value = "94.03"
$ws2.Range("F29").NumberFormat = "@"
$ws2.Range("F29").Value = "0.76"
$ws2.Range("G29").NumberFormat = "@"
$ws2.Range("G29").Value = "0.0027"
$ws2.Range("H29").Value = 7

# row 30: index 28 -> 014930 富国高端制造行业股票C
$cell = $ws2.Range("A30"); $styleSrc.Copy($cell); $cell.Value = 28
$ws2.Range("B30").NumberFormat = "@"
$ws2.Range("B30").Value = "014930"
$ws2.Range("C30").Value = "富国高端制造行业股票C"
$ws2.Range("D30").NumberFormat = "@"
$ws2.Range("D30").Value = "0.01"
$ws2.Range("E30").NumberFormat = "@"
$ws2.Range("E30").Value = "91.41"
$ws2.Range("F30").NumberFormat = "@"
$ws2.Range("F30").Value = "3.12"
$ws2.Range("G30").NumberFormat = "@"
$ws2.Range("G30").Value = "0.0003"
$ws2.Range("H30").Value = 9

